$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.409.69"
$ws.Range("E2").Value = "  +2.52%  "
$ws.Range("D3").Value = "2.109.77"
$ws.Range("E3").Value = "  +0.67%  "
$ws.Range("D4").Value = "1.007"
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "345.25"
$ws.Range("E6").Value = "  -0.08%  "
$ws.Range("D7").Value = "0.5234"
$ws.Range("E7").Value = "  +2.19%  "
$ws.Range("D8").Value = "0.4444"
$ws.Range("E8").Value = "  +1.09%  "
$ws.Range("D9").Value = "54.75"
$ws.Range("E9").Value = "  +2.34%  "
$ws.Range("D10").Value = "0.09380"
$ws.Range("E10").Value = "  +2.80%  "
$ws.Range("D11").Value = "1.174"
$ws.Range("D12").Value = "25.01"
$ws.Range("E12").Value = "  +0.82%  "
$ws.Range("D13").Value = "8.690"
$ws.Range("E13").Value = "  +6.04%  "
$ws.Range("D14").Value = "6.951"
$ws.Range("E14").Value = "  +3.17%  "
$ws.Range("D15").Value = "2.022.53"
$ws.Range("E15").Value = "  -3.52%  "
$ws.Range("D16").Value = "101.99"
$ws.Range("E16").Value = "  +2.40%  "
$ws.Range("D17").Value = "0.00001163"
$ws.Range("E17").Value = "  +1.59%  "
$ws.Range("D18").Value = "1.008"
$ws.Range("E18").Value = "  -0.03%  "
$ws.Range("D19").Value = "21.22"
$ws.Range("E19").Value = "  +0.65%  "
$ws.Range("D20").Value = "0.06726"
$ws.Range("E20").Value = "  +1.30%  "
$ws.Range("D21").Value = "6.342"
$ws.Range("E21").Value = "  +2.71%  "
$ws.Range("D22").Value = "1.006"
$ws.Range("E22").Value = "  -0.05%  "
$ws.Range("D23").Value = "30.442.90"
$ws.Range("E23").Value = "  +2.45%  "
$ws.Range("D24").Value = "12.66"
$ws.Range("E24").Value = "  +0.65%  "
$ws.Range("D25").Value = "2.298"
$ws.Range("E25").Value = "  -0.50%  "
$ws.Range("D26").Value = "22.07"
$ws.Range("E26").Value = "  +1.19%  "
$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D27").Value = "163.03"
$ws.Range("E27").Value = "  +0.30%  "
$ws.Range("B28").Value = "LidoDAOToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D28").Value = "2.536"
$ws.Range("E28").Value = "  +0.66%  "
$ws.Range("D29").Value = "134.25"
$ws.Range("E29").Value = "  +1.39%  "
$ws.Range("D30").Value = "1.155"
$ws.Range("E30").Value = "  +2.33%  "
$ws.Range("D31").Value = "1.739"
$ws.Range("E31").Value = "  +6.38%  "
$ws.Range("E32").Value = "  +1.20%  "
$ws.Range("D33").Value = "6.815"
$ws.Range("E33").Value = "  +12.77%  "
$ws.Range("D34").Value = "6.272"
$ws.Range("E34").Value = "  +1.87%  "
$ws.Range("E35").Value = "  -1.05%  "
$ws.Range("E36").Value = "  +1.35%  "
$ws.Range("D37").Value = "0.02624"
$ws.Range("E37").Value = "  +2.15%  "
$ws.Range("D38").Value = "0.06795"
$ws.Range("E38").Value = "  +1.92%  "
$ws.Range("D39").Value = "0.7058"
$ws.Range("E39").Value = "  +3.08%  "
$ws.Range("D40").Value = "1.357"
$ws.Range("E40").Value = "  +5.08%  "
$ws.Range("D41").Value = "12.58"
$ws.Range("E41").Value = "  +1.75%  "
$ws.Range("D42").Value = "0.2227"
$ws.Range("E42").Value = "  -0.13%  "
$ws.Range("D43").Value = "0.6860"
$ws.Range("E43").Value = "  +2.59%  "
$ws.Range("D44").Value = "14.52"
$ws.Range("E44").Value = "  +2.63%  "
$ws.Range("E45").Value = "  +3.20%  "
$ws.Range("D46").Value = "1.005"
$ws.Range("E46").Value = "  -0.02%  "
$ws.Range("D47").Value = "1.366"
$ws.Range("E47").Value = "  +17.53%  "
$ws.Range("E48").Value = "  +1.13%  "
$ws.Range("D49").Value = "0.00000000345"
$ws.Range("E49").Value = "  +3.16%  "
$ws.Range("B50").Value = "EOS"
$ws.Range("C50").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D50").Value = "1.221"
$ws.Range("E50").Value = "  +0.24%  "
$ws.Range("B51").Value = "ThetaToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D51").Value = "1.205"
$ws.Range("E51").Value = "  +8.95%  "
